$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.381.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.369.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.698"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.83"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.41%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.598"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +27.92%  "
$ws.Range("E10").Value = "  +5.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "31.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +17.22%  "
$ws.Range("E12").Value = "  +18.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.723.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "17.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.921"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.371.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.404.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "  +3.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "79.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "258.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.93%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  +4.17%  "
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.36%  "
$ws.Range("E27").Value = "  +1.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.130"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.136"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0763"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.50%  "
$ws.Range("E35").Value = "  +6.04%  "
$ws.Range("E36").Value = "  +6.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0277"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.06%  "
$ws.Range("E41").Value = "  +2.01%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("E43").Value = "  +19.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.101"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.95%  "
$ws.Range("E45").Value = "  +2.83%  "
$ws.Range("E46").Value = "  +3.81%  "
$ws.Range("E47").Value = "  +11.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("E49").Value = "  -1.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.472.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.33%  "
